$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill previously-blank cells in column A (rows 2-8) with 0
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 0
$ws.Range("A5").Value = 0
$ws.Range("A6").Value = 0
$ws.Range("A7").Value = 0
$ws.Range("A8").Value = 0

# Fill the other previously-blank cells with 0
$ws.Range("B4").Value = 0
$ws.Range("B6").Value = 0

# Overwrite C5 (previously the shared string "test_5") with the numeric 0
$ws.Range("C5").Value = 0

# Turn on AutoFilter over the used data range
[void]$ws.Range("A1:D8").AutoFilter()

# Register the hidden _FilterDatabase defined name that Excel creates
# for the sheet-level autofilter (localSheetId=0, hidden=1)
$filterRange = $ws.Range("A1:D8")
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $filterRange)
$filterName.Visible = $false

# Move the active selection, matching the end-state cursor position
[void]$ws.Range("M14").Select()
